# Updated cryptos list (price + 1h volume columns) to match the
# Tue Aug  1 14:45:52 UTC 2023 GitHub Actions refresh.
#
# D (Price) and E (Volume(1h)) are plain text cells in the source sheet.
# Values that look like pure numbers (e.g. "1.000", "0.9994") would be
# auto-converted to numeric cells by Excel's normal value parsing, so
# those are written with a leading apostrophe (quote-prefix) to force them
# to stay text, exactly like the original data. Values that already are not
# parseable as numbers (e.g. "28.857.05", the Volume(1h) percentages with
# surrounding spaces) are assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.857.05'
$ws.Range("E2").Value = '  -1.65%  '
# Row 3
$ws.Range("D3").Value = '1.833.43'
$ws.Range("E3").Value = '  -1.66%  '
# Row 4
$ws.Range("D4").Value = '''0.9994'
$ws.Range("E4").Value = '  -0.10%  '
# Row 5
$ws.Range("D5").Value = '''244.56'
$ws.Range("E5").Value = '  +0.51%  '
# Row 6
$ws.Range("D6").Value = '''0.6929'
$ws.Range("E6").Value = '  -1.11%  '
# Row 7
$ws.Range("E7").Value = '  -0.07%  '
# Row 8
$ws.Range("D8").Value = '''0.07684'
$ws.Range("E8").Value = '  -2.86%  '
# Row 9
$ws.Range("D9").Value = '''0.3048'
$ws.Range("E9").Value = '  -2.44%  '
# Row 10
$ws.Range("D10").Value = '''23.21'
$ws.Range("E10").Value = '  -4.74%  '
# Row 11
$ws.Range("D11").Value = '''0.07795'
$ws.Range("E11").Value = '  -0.07%  '
# Row 12
$ws.Range("D12").Value = '''93.03'
$ws.Range("E12").Value = '  +0.75%  '
# Row 13
$ws.Range("D13").Value = '1.834.95'
$ws.Range("E13").Value = '  -1.86%  '
# Row 14
$ws.Range("D14").Value = '''5.089'
$ws.Range("E14").Value = '  -1.08%  '
# Row 15
$ws.Range("D15").Value = '''0.6792'
$ws.Range("E15").Value = '  -2.70%  '
# Row 16
$ws.Range("D16").Value = '''6.446'
$ws.Range("E16").Value = '  -1.74%  '
# Row 17
$ws.Range("E17").Value = '  -3.27%  '
# Row 18
$ws.Range("D18").Value = '28.860.85'
$ws.Range("E18").Value = '  -1.74%  '
# Row 19
$ws.Range("D19").Value = '''242.45'
$ws.Range("E19").Value = '  -2.52%  '
# Row 20
$ws.Range("D20").Value = '2.074.13'
$ws.Range("E20").Value = '  -2.46%  '
# Row 21
$ws.Range("D21").Value = '''12.69'
$ws.Range("E21").Value = '  -2.30%  '
# Row 22
$ws.Range("D22").Value = '''1.000'
$ws.Range("E22").Value = '  +0.01%  '
# Row 23
$ws.Range("D23").Value = '''7.438'
$ws.Range("E23").Value = '  -1.97%  '
# Row 24
$ws.Range("D24").Value = '''0.9999'
# Row 25
$ws.Range("D25").Value = '''0.1484'
$ws.Range("E25").Value = '  -3.39%  '
# Row 26
$ws.Range("D26").Value = '''159.49'
$ws.Range("E26").Value = '  -0.78%  '
# Row 27
$ws.Range("D27").Value = '''8.769'
$ws.Range("E27").Value = '  -2.24%  '
# Row 28
$ws.Range("D28").Value = '''18.23'
$ws.Range("E28").Value = '  -2.74%  '
# Row 29
$ws.Range("D29").Value = '''1.540'
$ws.Range("E29").Value = '  -3.32%  '
# Row 30
$ws.Range("D30").Value = '''4.216'
$ws.Range("E30").Value = '  -1.94%  '
# Row 31
$ws.Range("D31").Value = '''4.155'
$ws.Range("E31").Value = '  -2.05%  '
# Row 32
$ws.Range("D32").Value = '''1.184'
$ws.Range("E32").Value = '  -1.76%  '
# Row 33
$ws.Range("D33").Value = '''0.05094'
$ws.Range("E33").Value = '  -2.91%  '
# Row 34
$ws.Range("D34").Value = '''0.7742'
$ws.Range("E34").Value = '  +2.37%  '
# Row 35
$ws.Range("D35").Value = '''1.855'
$ws.Range("E35").Value = '  -1.59%  '
# Row 36
$ws.Range("D36").Value = '''1.141'
$ws.Range("E36").Value = '  -3.29%  '
# Row 37
$ws.Range("D37").Value = '''2.692'
$ws.Range("E37").Value = '  -0.28%  '
# Row 38
$ws.Range("D38").Value = '''0.01848'
$ws.Range("E38").Value = '  -0.95%  '
# Row 39
$ws.Range("D39").Value = '1.242.55'
$ws.Range("E39").Value = '  -2.66%  '
# Row 40
$ws.Range("D40").Value = '''2.698'
$ws.Range("E40").Value = '  -1.83%  '
# Row 41
$ws.Range("D41").Value = '''0.9500'
$ws.Range("E41").Value = '  +5.55%  '
# Row 42
$ws.Range("D42").Value = '''107.75'
$ws.Range("E42").Value = '  -1.86%  '
# Row 43
$ws.Range("D43").Value = '''5.952'
$ws.Range("E43").Value = '  -0.55%  '
# Row 44
$ws.Range("E44").Value = '  +0.01%  '
# Row 45
$ws.Range("D45").Value = '''9.610'
$ws.Range("E45").Value = '  +0.05%  '
# Row 46
$ws.Range("D46").Value = '1.975.70'
$ws.Range("E46").Value = '  -2.28%  '
# Row 47
$ws.Range("D47").Value = '''0.5157'
$ws.Range("E47").Value = '  -0.34%  '
# Row 48
$ws.Range("D48").Value = '''63.83'
$ws.Range("E48").Value = '  -9.04%  '
# Row 49
$ws.Range("D49").Value = '''1.742'
$ws.Range("E49").Value = '  -2.74%  '
# Row 50
$ws.Range("D50").Value = '''0.00000000116'
$ws.Range("E50").Value = '  -7.87%  '
# Row 51
$ws.Range("D51").Value = '''6.919'
$ws.Range("E51").Value = '  -1.42%  '
